$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "username/user/password" list with the new Amazon
# exercise items, and drop the now-unused second (B) column entirely.
$ws.Range("A1").Value = "Baquetas"
$ws.Range("A2").Value = "Nintendo Switch"
$ws.Range("A3").Value = "Audifonos Inalambricos"

$ws.Range("B1:B3").Delete()
